$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '28.576.96' },
    @{ Cell = 'E2'; Value = '  -3.23%  ' },
    @{ Cell = 'D3'; Value = '1.850.75' },
    @{ Cell = 'E3'; Value = '  -3.60%  ' },
    @{ Cell = 'D4'; Value = '1.002' },
    @{ Cell = 'E4'; Value = '  -1.09%  ' },
    @{ Cell = 'D5'; Value = '335.53' },
    @{ Cell = 'E5'; Value = '  +2.92%  ' },
    @{ Cell = 'E6'; Value = '  -0.95%  ' },
    @{ Cell = 'D7'; Value = '0.4666' },
    @{ Cell = 'E7'; Value = '  -3.13%  ' },
    @{ Cell = 'D8'; Value = '0.3918' },
    @{ Cell = 'E8'; Value = '  -3.46%  ' },
    @{ Cell = 'D9'; Value = '46.52' },
    @{ Cell = 'E9'; Value = '  -2.90%  ' },
    @{ Cell = 'D10'; Value = '0.07884' },
    @{ Cell = 'E10'; Value = '  -4.14%  ' },
    @{ Cell = 'D11'; Value = '0.9850' },
    @{ Cell = 'E11'; Value = '  -2.46%  ' },
    @{ Cell = 'D12'; Value = '22.19' },
    @{ Cell = 'E12'; Value = '  -5.27%  ' },
    @{ Cell = 'D13'; Value = '1.962.92' },
    @{ Cell = 'E13'; Value = '  +1.00%  ' },
    @{ Cell = 'D14'; Value = '5.854' },
    @{ Cell = 'E14'; Value = '  -3.39%  ' },
    @{ Cell = 'D15'; Value = '7.026' },
    @{ Cell = 'E15'; Value = '  -3.06%  ' },
    @{ Cell = 'D16'; Value = '0.06832' },
    @{ Cell = 'E16'; Value = '  -0.52%  ' },
    @{ Cell = 'D17'; Value = '87.67' },
    @{ Cell = 'E17'; Value = '  -4.26%  ' },
    @{ Cell = 'E18'; Value = '  -1.11%  ' },
    @{ Cell = 'E19'; Value = '  -2.83%  ' },
    @{ Cell = 'E20'; Value = '  -3.01%  ' },
    @{ Cell = 'D21'; Value = '1.002' },
    @{ Cell = 'E21'; Value = '  -0.97%  ' },
    @{ Cell = 'D22'; Value = '28.589.57' },
    @{ Cell = 'E22'; Value = '  -3.20%  ' },
    @{ Cell = 'D23'; Value = '5.410' },
    @{ Cell = 'E23'; Value = '  -4.64%  ' },
    @{ Cell = 'D24'; Value = '11.27' },
    @{ Cell = 'E24'; Value = '  -4.91%  ' },
    @{ Cell = 'D25'; Value = '2.122' },
    @{ Cell = 'E25'; Value = '  -3.29%  ' },
    @{ Cell = 'D26'; Value = '2.117.86' },
    @{ Cell = 'E26'; Value = '  -2.24%  ' },
    @{ Cell = 'D27'; Value = '153.62' },
    @{ Cell = 'E27'; Value = '  -1.56%  ' },
    @{ Cell = 'D28'; Value = '6.165' },
    @{ Cell = 'E28'; Value = '  -5.95%  ' },
    @{ Cell = 'D30'; Value = '2.024' },
    @{ Cell = 'E30'; Value = '  -3.66%  ' },
    @{ Cell = 'D31'; Value = '117.73' },
    @{ Cell = 'E31'; Value = '  -2.42%  ' },
    @{ Cell = 'D32'; Value = '0.9761' },
    @{ Cell = 'E32'; Value = '  -4.23%  ' },
    @{ Cell = 'E33'; Value = '  -2.02%  ' },
    @{ Cell = 'D34'; Value = '5.380' },
    @{ Cell = 'E34'; Value = '  -4.28%  ' },
    @{ Cell = 'E35'; Value = '  -1.62%  ' },
    @{ Cell = 'D36'; Value = '1.349' },
    @{ Cell = 'E36'; Value = '  -1.88%  ' },
    @{ Cell = 'D37'; Value = '0.06116' },
    @{ Cell = 'E37'; Value = '  -3.06%  ' },
    @{ Cell = 'E38'; Value = '  -3.97%  ' },
    @{ Cell = 'D39'; Value = '1.162' },
    @{ Cell = 'E39'; Value = '  -1.92%  ' },
    @{ Cell = 'B40'; Value = 'Frax' },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax' },
    @{ Cell = 'D40'; Value = '1.002' },
    @{ Cell = 'E40'; Value = '  -0.89%  ' },
    @{ Cell = 'B41'; Value = 'TheSandbox' },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand' },
    @{ Cell = 'D41'; Value = '0.5705' },
    @{ Cell = 'E41'; Value = '  -4.00%  ' },
    @{ Cell = 'B42'; Value = 'FraxShare' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' },
    @{ Cell = 'D42'; Value = '7.589' },
    @{ Cell = 'E42'; Value = '  -4.03%  ' },
    @{ Cell = 'B43'; Value = 'Aptos' },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' },
    @{ Cell = 'D43'; Value = '10.12' },
    @{ Cell = 'E43'; Value = '  -5.68%  ' },
    @{ Cell = 'B44'; Value = 'Algorand' },
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' },
    @{ Cell = 'D44'; Value = '0.1792' },
    @{ Cell = 'E44'; Value = '  -3.02%  ' },
    @{ Cell = 'B45'; Value = 'RenderToken' },
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Cell = 'D45'; Value = '2.378' },
    @{ Cell = 'E45'; Value = '  -3.36%  ' },
    @{ Cell = 'B46'; Value = 'WEMIXToken' },
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' },
    @{ Cell = 'D46'; Value = '1.260' },
    @{ Cell = 'E46'; Value = '  -1.66%  ' },
    @{ Cell = 'B47'; Value = 'EnergySwap' },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Cell = 'D47'; Value = '11.85' },
    @{ Cell = 'E47'; Value = '  -4.15%  ' },
    @{ Cell = 'B48'; Value = 'Decentraland' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' },
    @{ Cell = 'D48'; Value = '0.5394' },
    @{ Cell = 'E48'; Value = '  -3.09%  ' },
    @{ Cell = 'B49'; Value = 'Cronos' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' },
    @{ Cell = 'D49'; Value = '0.07160' },
    @{ Cell = 'E49'; Value = '  -4.15%  ' },
    @{ Cell = 'B50'; Value = 'NEARProtocol' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' },
    @{ Cell = 'D50'; Value = '1.907' },
    @{ Cell = 'E50'; Value = '  -1.86%  ' },
    @{ Cell = 'B51'; Value = 'Quant' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt' },
    @{ Cell = 'D51'; Value = '113.81' },
    @{ Cell = 'E51'; Value = '  -4.07%  ' }
)

foreach ($u in $updates) {
    $cellRef = $u.Cell
    $val = $u.Value
    $isNumeric = $val -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$'
    if ($isNumeric) {
        $ws.Range($cellRef).NumberFormat = '@'
        $ws.Range($cellRef).Value = $val
        $ws.Range($cellRef).Style = 'Normal'
    } else {
        $ws.Range($cellRef).Value = $val
    }
}
